$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "SamplePortion"
$ws.Range("K1").Value = "SamplePortionUnit"
$ws.Range("L1").Value = "Solvent"
$ws.Range("M1").Value = "TemperatureCycle"
$ws.Range("N1").Value = "BucketType"
$ws.Range("O1").Value = "DryingTemperature"
$ws.Range("P1").Value = "ResultUnit"
$ws.Range("Q1").Value = "Result"
$ws.Range("R1").Value = "Comment"

$ws.Range("J2").Value = "# Prise d'essai"
$ws.Range("K2").Value = "# Unité de mesure de la prise d’essai"
$ws.Range("L2").Value = "#Solvant"
$ws.Range("M2").Value = "#CycleDeTemperature"
$ws.Range("N2").Value = "#TypeDeGodet"
$ws.Range("O2").Value = "#TemperatureDeChauffage"
$ws.Range("P2").Value = "# Unité du résultat"
$ws.Range("Q2").Value = "# Résultat"
$ws.Range("R2").Value = "# Commentaire"

$ws.Range("J3").Value = "#float"
$ws.Range("K3").Value = "#string"
$ws.Range("L3").Value = "#string"
$ws.Range("M3").Value = "#string"
$ws.Range("N3").Value = "#string"
$ws.Range("O3").Value = "#float,`n  unit:celsius"
$ws.Range("P3").Value = "#string"
$ws.Range("Q3").Value = "#float"
$ws.Range("R3").Value = "#string"

$ws.Range("J4").Value = "# format: nombre décimal, ne pas spécifier d'unité"
$ws.Range("K4").Value = "# format: texte"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = "# format: texte"
$ws.Range("Q4").Value = "# format: nombre décimal ou NA"
$ws.Range("R4").Value = "# format: texte libre"

$ws.Range("J5").Value = "# ex: 2.5"
$ws.Range("K5").Value = "# ex: mg"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = "# ex: mg/ml"
$ws.Range("Q5").Value = "# 409.935 ou NA"
$ws.Range("R5").Value = ""
